$d = $word.ActiveDocument

# --- Change 1: "Link a tarea" paragraph -- drop en-US lang formatting,
#     restructure proofErr around "Link" (gramStart/End) instead of "tarea" (spellStart/End)
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4.InsertXML('<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>Link</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> a tarea: </w:t></w:r><w:hyperlink r:id="rId4" w:history="1"><w:r><w:t>https://github.com/jurrutiag/cc3002-tarea1</w:t></w:r></w:hyperlink></w:p>')

# Re-apply the Hyperlink character style (InsertXML drops rStyle refs)
$hl = $d.Hyperlinks.Item(1)
$hl.Range.Style = "Hipervnculo"

# --- Change 2: merge the "propiedades:" paragraph with the following
#     bookmark+drawing paragraph, dropping the _GoBack bookmark from here
$p6 = $d.Paragraphs.Item(6)
$p7 = $d.Paragraphs.Item(7)
$rngMerge = $d.Range($p6.Range.Start, $p7.Range.End)
$rngMerge.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:r><w:t>Diagrama UML mostrando solo métodos y</w:t></w:r><w:r><w:t xml:space="preserve"> propiedades:</w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="4E13AB03" wp14:editId="639471D9"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>-332105</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>285115</wp:posOffset></wp:positionV><wp:extent cx="9907905" cy="5965190"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapThrough wrapText="bothSides"><wp:wrapPolygon edited="0"><wp:start x="0" y="0"/><wp:lineTo x="0" y="21522"/><wp:lineTo x="21554" y="21522"/><wp:lineTo x="21554" y="0"/><wp:lineTo x="0" y="0"/></wp:wrapPolygon></wp:wrapThrough><wp:docPr id="3" name="Imagen 3"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 5"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId6" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="9907905" cy="5965190"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r></w:p>')

# --- Change 3: replace the trailing empty paragraph with the new
#     "Package Pokemons..." text, underline pPr mark, and _GoBack bookmark
$pLast = $d.Paragraphs.Last
$rngLast = $pLast.Range
$rngLast.InsertXML('<w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Package</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Pokemons</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Package</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Attacks</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> no están expandidos de</w:t></w:r><w:r><w:t xml:space="preserve">bido a que son ataques y </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pokemones</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> específicos, creados para facilitar el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>testing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. No tienen métodos ni parámetros extra en comparación con la superclase de cada uno.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
